$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 20 (all subsequent rows,
# previously 20-79, shift down to 21-80). Insert a blank row above the
# current row 20 first so everything below shifts down, then fill it in.
$ws.Rows.Item(20).Insert()

$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(20, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(20, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(20, 4).Value = 44648
$ws.Cells.Item(20, 5).Value = 15
$ws.Cells.Item(20, 6).Value = 100112038
$ws.Cells.Item(20, 7).Value = "Cebollín baby"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 300
$ws.Cells.Item(20, 11).Value = 1300
$ws.Cells.Item(20, 12).Value = 1500
$ws.Cells.Item(20, 13).Value = 1400
$ws.Cells.Item(20, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(20, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(20, 16).Value = 700
$ws.Cells.Item(20, 17).Value = 2
$ws.Cells.Item(20, 18).Value = "Hortaliza"
